$wb = $excel.ActiveWorkbook

# --- WECC sheet (sheet1): zero out the "Fixed cost" column (C2:C30) ---
$wsWecc = $wb.Worksheets.Item("WECC")
$wsWecc.Range("C2:C30").Value = 0

# --- NPCC sheet (sheet2): zero out the "Fixed cost" column (C2:C49) ---
$wsNpcc = $wb.Worksheets.Item("NPCC")
$wsNpcc.Range("C2:C49").Value = 0

# Widen column B on NPCC so the longer "Incremental cost" values are visible
$wsNpcc.Columns.Item(2).ColumnWidth = 35.33

# Restore the active selections left by the editor on each sheet
$wsWecc.Range("F13").Select() | Out-Null
$wsNpcc.Range("C37:C49").Select() | Out-Null

# Make WECC the active (tab-selected) sheet, matching the saved workbook state
$wsWecc.Activate() | Out-Null
